# Auto-generated edit script: clean up footnote markers ([1],[2],[3],[4],[5],[5, 6])
# from Vaccine names, merge multi-line cell text into single-line text (joined with
# spaces), and de-duplicate shared strings that become identical after cleanup.
$wb = $excel.ActiveWorkbook

# --- Sheet 1: Pediatric Vaccine  ---
$ws = $wb.Worksheets.Item(1)
$ws.Range('A2').Value = 'DTaP '
$ws.Range('A3').Value = 'DTaP '
$ws.Range('A4').Value = 'DTaP '
$ws.Range('A5').Value = 'DTaP '
$ws.Range('A6').Value = 'DTaP-IPV '
$ws.Range('A7').Value = 'DTaP-IPV '
$ws.Range('A8').Value = 'DTaP-Hep B-IPV '
$ws.Range('A9').Value = 'DTaP-IP-HI '
$ws.Range('A10').Value = 'e-IPV '
$ws.Range('A11').Value = 'Hepatitis B-Hib '
$ws.Range('A12').Value = 'Hepatitis A Pediatric '
$ws.Range('A13').Value = 'Hepatitis A Pediatric '
$ws.Range('A14').Value = 'Hepatitis A Pediatric '
$ws.Range('A15').Value = 'Hepatitis A-Hepatitis B 18 only '
$ws.Range('A16').Value = 'Hepatitis A-Hepatitis B 18 only '
$ws.Range('A17').Value = 'Hepatitis B  Pediatric/Adolescent'
$ws.Range('A18').Value = 'Hepatitis B  Pediatric/Adolescent'
$ws.Range('A19').Value = 'Hepatitis B  Pediatric/Adolescent'
$ws.Range('B19').Value = 'Recombivax HB'
$ws.Range('A20').Value = 'Hib '
$ws.Range('A21').Value = 'Hib '
$ws.Range('A22').Value = 'Hib '
$ws.Range('A23').Value = 'HPV - Quadrivalent Human Papillomavirus Types 6, 11, 16 and 18 Recombinant '
$ws.Range('A24').Value = 'HPV -Bivalent Human Papillomavirus Types 16 and 18 '
$ws.Range('A25').Value = 'Meningococcal Conjugate (Groups A, C, Y and W-135) '
$ws.Range('A26').Value = 'Meningococcal Conjugate (Groups A, C, Y and W-135) '
$ws.Range('A27').Value = 'Measles, Mumps and Rubella (MMR) '
$ws.Range('A28').Value = 'Pneumococcal 13-valent  (Pediatric)'
$ws.Range('A30').Value = 'Rotavirus, Live, Oral, Pentavalent '
$ws.Range('A31').Value = 'Rotavirus, Live, Oral, Oral '
$ws.Range('A32').Value = 'Tetanus  Diphtheria Toxoids '
$ws.Range('D32').Value = '10 pack - 1 dose syringes No Needle'
$ws.Range('A33').Value = 'Tetanus  Diphtheria Toxoids '
$ws.Range('A34').Value = 'Tetanus  Diphtheria Toxoids '
$ws.Range('A35').Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range('A36').Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range('D36').Value = '10 pack - 1 dose TL syringes, No Needle'
$ws.Range('A37').Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range('A38').Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range('A39').Value = 'Varicella '

# --- Sheet 2: Adult Vaccine  ---
$ws = $wb.Worksheets.Item(2)
$ws.Range('A2').Value = 'Hepatitis A Adult '
$ws.Range('A3').Value = 'Hepatitis A Adult '
$ws.Range('A4').Value = 'Hepatitis A-Hepatitis B Adult '
$ws.Range('A5').Value = 'Hepatitis A-Hepatitis B Adult '
$ws.Range('A6').Value = 'Hepatitis B-Adult '
$ws.Range('A7').Value = 'Hepatitis B-Adult '
$ws.Range('A8').Value = 'HPV -Quadrivalent Human Papillomavirus Types 6, 11, 16 and 18 Recombinant Adult '
$ws.Range('A9').Value = 'HPV-Human Papillomavirus Bivalent Types 16 and 18 '
$ws.Range('A10').Value = 'HPV-Human Papillomavirus Bivalent Types 16 and 18 '
$ws.Range('A11').Value = 'Measles, Mumps,  Rubella-Adult '
$ws.Range('A14').Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range('A15').Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range('A16').Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range('A17').Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range('A18').Value = 'Varicella-Adult '
$ws.Range('A21').Value = 'Tetanus and Diphtheria Toxoids '
$ws.Range('A22').Value = 'Meningococcal Conjugate (Groups A, C, H and W-135) '

# --- Sheet 3: Pediatric Influenza Vaccine  ---
$ws = $wb.Worksheets.Item(3)
$ws.Range('A2').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('A3').Value = 'Influenza  (Age 6-35 months)'
$ws.Range('B3').Value = 'Fluzone Pediatric dose No Preservative'
$ws.Range('A4').Value = 'Influenza  (Age 36 months and older)'
$ws.Range('B4').Value = 'Fluzone No-Preservative'
$ws.Range('A5').Value = 'Influenza  (Age 36 months and older)'
$ws.Range('B5').Value = 'Fluzone No-Preservative'
$ws.Range('A6').Value = 'Influenza  (Age 36 months and older)'
$ws.Range('B6').Value = 'Fluarix Preservative Free'
$ws.Range('D6').Value = '10 pack- 1 dose TipLok syringe'
$ws.Range('A7').Value = 'Influenza  (Age 4 years and older)'
$ws.Range('A8').Value = 'Influenza  (Age 4 years and older)'
$ws.Range('B8').Value = 'Fluvirin Preservative-free'
$ws.Range('A9').Value = 'Influenza  Live, Intranasal (Age 2-49 years)'
$ws.Range('B9').Value = 'FluMist No Preservative'
$ws.Range('A10').Value = 'Influenza  (Age 9 years and older)'
$ws.Range('B10').Value = 'Afluria No Preservative'
$ws.Range('D10').Value = '10 pack-1 dose syringe'
$ws.Range('H10').Value = 'Merck (CSL product)'
$ws.Range('A11').Value = 'Influenza  (Age 9 years and older)'
$ws.Range('H11').Value = 'Merck (CSL product)'

# --- Sheet 4: Adult Influenza Vaccine  ---
$ws = $wb.Worksheets.Item(4)
$ws.Range('A2').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('A3').Value = 'Influenza  (age 36 months and older)'
$ws.Range('B3').Value = 'Fluzone No Preservative'
$ws.Range('A4').Value = 'Influenza  (age 36 months and older)'
$ws.Range('B4').Value = 'Fluzone No Preservative'
$ws.Range('A5').Value = 'Influenza  (Age 4 years and older)'
$ws.Range('A6').Value = 'Influenza  (Age 4 years and older)'
$ws.Range('B6').Value = 'Fluvirin Preservative-free'
$ws.Range('A7').Value = 'Influenza  (age 36 months and older)'
$ws.Range('A8').Value = 'Influenza  (18 years and older)'
$ws.Range('A9').Value = 'Influenza  Live, Intranasal (Age 2-49 years)'
$ws.Range('B9').Value = 'FluMist  No Preservative'
$ws.Range('A10').Value = 'Influenza  (Age 9 years and older)'
$ws.Range('B10').Value = 'Afluria No Preservative'
$ws.Range('D10').Value = '10 pack-1 dose syringe'
$ws.Range('H10').Value = 'Merck (CSL product)'
$ws.Range('A11').Value = 'Influenza  (Age 9 years and older)'
$ws.Range('H11').Value = 'Merck (CSL product)'

